$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("category")

$ws.Range("A7").Value = "Feria"
$ws.Range("B7").Value = "Expense"

$ws.Range("A8").Value = "Mercado"
$ws.Range("B8").Value = "Expense"

$ws.Range("A9").Value = "Netflix"
$ws.Range("B9").Value = "Expense"
